$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.534624457359314
$ws.Range("B1").Value = 5.458083629608154
$ws.Range("C1").Value = 3.621304035186768
$ws.Range("D1").Value = 0.9503316283226013
$ws.Range("E1").Value = 0.6074384450912476
